# Update countries & provincias Spain
# Refresh COVID country statistics table and the "last updated" banner.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Datos actualizados ..." banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 12 de Agosto de 2020 a las 18:28"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 5322171
$ws.Range("C4").Value = 16214
$ws.Range("D4").Value = 2758628
$ws.Range("E4").Value = 2395287
$ws.Range("G4").Value = 507
$ws.Range("H4").Value = 168256

# India (row 6)
$ws.Range("B6").Value = 2372318
$ws.Range("C6").Value = 43913
$ws.Range("D6").Value = 1673885
$ws.Range("E6").Value = 651716
$ws.Range("G6").Value = 529
$ws.Range("H6").Value = 46717

# Reino Unido (row 15)
$ws.Range("B15").Value = 313798
$ws.Range("C15").Value = 1009
$ws.Range("G15").Value = 77
$ws.Range("H15").Value = 46706

# Italia (row 20)
$ws.Range("B20").Value = 251713
$ws.Range("C20").Value = 476
$ws.Range("D20").Value = 202697
$ws.Range("E20").Value = 13791
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = 35225

# Canada (row 27)
$ws.Range("B27").Value = 120617
$ws.Range("C27").Value = 196
$ws.Range("D27").Value = 107015
$ws.Range("E27").Value = 4598
$ws.Range("G27").Value = 13
$ws.Range("H27").Value = 9004

# Republica Dominicana (row 38)
$ws.Range("B38").Value = 82224
$ws.Range("C38").Value = 1130
$ws.Range("D38").Value = 47095
$ws.Range("E38").Value = 33758
$ws.Range("G38").Value = 25
$ws.Range("H38").Value = 1371

# Singapur (row 47)
$ws.Range("D47").Value = 50520
$ws.Range("E47").Value = 4848

# Kenia/Venezuela swap rank: row 65 becomes Kenia (updated numbers),
# row 66 becomes Venezuela (its previous, unchanged numbers)
$ws.Range("A65").Value = "Kenia"
$ws.Range("B65").Value = 28104
$ws.Range("C65").Value = 679
$ws.Range("D65").Value = 14610
$ws.Range("E65").Value = 13038
$ws.Range("G65").Value = 18
$ws.Range("H65").Value = 456

$ws.Range("A66").Value = "Venezuela"
$ws.Range("B66").Value = 27938
$ws.Range("C66").Value = 0
$ws.Range("D66").Value = 19706
$ws.Range("E66").Value = 7994
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 238

# Chequia (row 74)
$ws.Range("B74").Value = 18984
$ws.Range("C74").Value = 201
$ws.Range("D74").Value = 13407
$ws.Range("E74").Value = 5186

# Hong Kong (row 111)
$ws.Range("D111").Value = 3189
$ws.Range("E111").Value = 992

# Jordania (row 145)
$ws.Range("B145").Value = 1303
$ws.Range("C145").Value = 20
$ws.Range("D145").Value = 1215
$ws.Range("E145").Value = 77

# Islas Malvinas/Montserrat swap rank: row 213 becomes Islas Malvinas,
# row 214 becomes Montserrat (each keeping its own previous numbers)
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0

$ws.Range("A214").Value = "Montserrat"
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
